$d = $word.ActiveDocument

# Locate the last-list-item paragraph, currently split across two runs
# ("Sta" / bookmark _GoBack / "vljanje IS u upotrebu") that together read
# "Stavljanje IS u upotrebu".
$target = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Stavljanje IS u upotrebu") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Stavljanje IS u upotrebu' paragraph"
}

$rng = $target.Range

# Replace that single paragraph with three paragraphs:
#  1) "Stavljanje IS u upotrebu" as one unified run (no more split / no list
#     numbering - plain body paragraph), sz 28.
#  2) An empty paragraph carrying the _GoBack bookmark, sz 28.
#  3) A new struck-through paragraph: "Sonarr, Raddar, plex," + " povezano je kroz API".
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>Stavljanje IS u upotrebu</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
<w:p><w:pPr><w:rPr><w:strike/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:sz w:val="28"/></w:rPr><w:t>Sonarr, Raddar, plex,</w:t></w:r><w:r><w:rPr><w:strike/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> povezano je kroz API</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$rng.InsertXML($xml)

# InsertXML replaces the range's contents but leaves the original paragraph
# mark behind as a trailing empty paragraph; collapse it away.
$countAfterInsert = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($countAfterInsert - 1)
$last = $d.Paragraphs.Item($countAfterInsert)
$cleanup = $d.Range($secondLast.Range.End - 1, $last.Range.End)
[void]$cleanup.Delete()
